$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 199, pushing existing rows 199-211 down to 200-212.
$ws.Range("A199").EntireRow.Insert()

# Populate the new row 199 with the weekly data point.
$ws.Range("A199").Value = 11
$ws.Range("B199").Value = "Vega Monumental Concepción"
$ws.Range("C199").Value = "Bíobío"
$ws.Range("D199").Value = 44826
$ws.Range("E199").Value = 8
$ws.Range("F199").Value = 100112003
$ws.Range("G199").Value = "Ajo"
$ws.Range("H199").Value = "Chino"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 400
$ws.Range("K199").Value = 17000
$ws.Range("L199").Value = 18000
$ws.Range("M199").Value = 17500
$ws.Range("N199").Value = "$/caja 10 kilos"
$ws.Range("O199").Value = "China"
$ws.Range("P199").Value = 1750
$ws.Range("Q199").Value = 10
$ws.Range("R199").Value = "Hortaliza"

# Match the date-number style used by the rest of column D.
$ws.Range("D199").NumberFormat = $ws.Range("D200").NumberFormat
